$wb = $excel.ActiveWorkbook

# Insert the new "Learning" sheet right after the last existing sheet ("Links")
$linksSheet = $wb.Worksheets.Item("Links")
$newSheet = $wb.Worksheets.Add($null, $linksSheet)
$newSheet.Name = "Learning"

# Populate the learning notes in column A
$newSheet.Range("A1").Value = "implement redux"
$newSheet.Range("A2").Value = "node backend to link with redux"
$newSheet.Range("A3").Value = "authentication and generating tokens"
$newSheet.Range("A4").Value = "responsive design"

# Match column A width used on the new sheet
$newSheet.Columns.Item(1).ColumnWidth = 28

# The newly added sheet becomes active; set its selection to match the target state
$newSheet.Range("A5").Select()
